# Printing all the graphs in one graph, calculating how long each child time travel
#
# The sheet stores every value as text (inlineStr in the original OOXML), so
# whenever a replacement value looks like a number Excel would normally
# coerce it to a numeric cell. Force "Text" number format first so the
# value committed to the cell is the literal string we want.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

function Set-TextCell {
    param($sheet, $row, $col, $value)
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# nChildren bumped from 7 to 8
Set-TextCell $ws 4 2 "8"

# Each child's row (car seat, first name, last name, coords, contact,
# pickup time, duration) was recalculated - this also adds a brand new
# 8th child (Ron Cohen) and shifts the trailing school/cost/time block
# down by one row.
$children = @(
    @{ Row = 6;  A = "0"; B = "14"; C = "Lorinda" + " " + $nbsp; D = "Tyron" + " " + $nbsp;     E = "3.37,7.62";  F = "Teresa(grandmother): 0558587699"; G = "7:00:00"; H = "42.0" },
    @{ Row = 7;  A = "1"; B = "16"; C = "Collette" + " " + $nbsp; D = "Billi" + " " + $nbsp;     E = "6.19,7.48";  F = "Elias(mother): 0578741979";       G = "7:05:00"; H = "37.0" },
    @{ Row = 8;  A = "2"; B = "18"; C = "Kandis" + " " + $nbsp;  D = "Zulma" + " " + $nbsp;      E = "9.32,9.44";  F = "Kylie(mother): 0575413269";       G = "7:10:00"; H = "32.0" },
    @{ Row = 9;  A = "3"; B = "11"; C = "Randolph" + " " + $nbsp; D = "Bridgette" + " " + $nbsp; E = "9.63,4.02";  F = "Lenny(father): 0505536740";       G = "7:17:00"; H = "25.0" },
    @{ Row = 10; A = "4"; B = "15"; C = "Nubia" + " " + $nbsp;   D = "Royce" + " " + $nbsp;      E = "9.71,4.53";  F = "Augustus(father): 0517389040";    G = "7:18:00"; H = "24.0" },
    @{ Row = 11; A = "5"; B = "10"; C = "Demetra" + " " + $nbsp; D = "Francene" + " " + $nbsp;   E = "8.81,2.38";  F = "Dorian(mother): 0534328089";      G = "7:21:00"; H = "21.0" },
    @{ Row = 12; A = "6"; B = "20"; C = "Ron";                    D = "Cohen";                    E = "9.6,-1.85";  F = "Bernardine(mother): 0576270618";  G = "7:28:00"; H = "14.0" },
    @{ Row = 13; A = "7"; B = "19"; C = "Jeanine" + " " + $nbsp; D = "Janee" + " " + $nbsp;      E = "6.54,0.52";  F = "Teresa(mother): 0517627420";      G = "7:34:00"; H = "8.0" }
)

foreach ($child in $children) {
    $r = $child.Row
    Set-TextCell $ws $r 1 $child.A
    Set-TextCell $ws $r 2 $child.B
    Set-TextCell $ws $r 3 $child.C
    Set-TextCell $ws $r 4 $child.D
    Set-TextCell $ws $r 5 $child.E
    Set-TextCell $ws $r 6 $child.F
    Set-TextCell $ws $r 7 $child.G
    Set-TextCell $ws $r 8 $child.H
}

# The trailing "school" / "cost" / "time" summary rows shift from
# 13/14/15 down to 14/15/16; the school row also gets a new total time
# in column G, and the overall time total becomes 42.0.
Set-TextCell $ws 14 1 "school"
Set-TextCell $ws 14 2 "3"
Set-TextCell $ws 14 3 "Ironiah"
Set-TextCell $ws 14 4 "mySchool"
Set-TextCell $ws 14 5 "0,0"
Set-TextCell $ws 14 6 "Shir(secretary): 0523345098"
Set-TextCell $ws 14 7 "7:42:00"

Set-TextCell $ws 15 1 "cost"
Set-TextCell $ws 15 2 "25"

Set-TextCell $ws 16 1 "time"
Set-TextCell $ws 16 2 "42.0"
